$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.314.51"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "3.573.49"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E5").Value = "  +2.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.41"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("D7").Value = "3.562.62"
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  +8.01%  "
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.56"
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "4.140.62"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("E16").Value = "  -1.17%  "
$ws.Range("D17").Value = "70.330.27"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").Value = "3.557.20"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "544.12"
$ws.Range("E21").Value = "  +10.84%  "
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.99"
$ws.Range("E23").Value = "  -7.60%  "
$ws.Range("E24").Value = "  +8.13%  "
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "96.01"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("E29").Value = "  -1.63%  "
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("E31").Value = "  -3.70%  "
$ws.Range("E32").Value = "  +3.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "65.23"
$ws.Range("E33").Value = "  -2.38%  "
$ws.Range("E34").Value = "  -1.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "555.29"
$ws.Range("E35").Value = "  -3.10%  "
$ws.Range("E36").Value = "  +3.87%  "
$ws.Range("E37").Value = "  +5.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.61"
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("E40").Value = "  -4.38%  "
$ws.Range("E41").Value = "  -4.16%  "
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("D43").Value = "3.365.48"
$ws.Range("E44").Value = "  -7.19%  "
$ws.Range("E45").Value = "  +3.27%  "
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("E47").Value = "  +1.59%  "
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.15"
$ws.Range("E49").Value = "  -5.45%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "137.43"
$ws.Range("E51").Value = "  +1.30%  "
